$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct district names to the official names from the website
$ws.Range("G7").Value = "Koppal"
$ws.Range("G8").Value = "Koppal"
$ws.Range("G34").Value = "Koppal"

# Remove the stray empty cell at F3 (it had no content and should not exist)
$ws.Range("F3").ClearContents()
